$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.000087185718439719
$ws.Range("D2").Value = 0.00000272839545435491
$ws.Range("E2").Value = 0.000019293960678235

$ws.Range("B3").Value = 2.69419413174802
$ws.Range("D3").Value = 20.9547780853451

$ws.Range("D4").Value = 23.7235235829361

$ws.Range("B5").Value = 44.3315961858685
$ws.Range("D5").Value = 265.838285402445
